$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 149
$ws.Range("D149").Value = 44516
$ws.Range("M149").Value = 36
$ws.Range("N149").Value = 7000
$ws.Range("O149").Value = 7000
$ws.Range("P149").Value = 7000
$ws.Range("R149").Value = 'Provincia de Quillota'
$ws.Range("S149").Value = 3500

# Row 150
$ws.Range("D150").Value = 44516
$ws.Range("L150").Value = 'Segunda'
$ws.Range("M150").Value = 38
$ws.Range("N150").Value = 5000
$ws.Range("O150").Value = 5000
$ws.Range("P150").Value = 5000
$ws.Range("S150").Value = 2500

# Row 151
$ws.Range("D151").Value = 44217
$ws.Range("M151").Value = 50
$ws.Range("N151").Value = 4000
$ws.Range("P151").Value = 4000
$ws.Range("S151").Value = 2000

# Row 152
$ws.Range("D152").Value = 44509
$ws.Range("M152").Value = 35
$ws.Range("N152").Value = 10000
$ws.Range("O152").Value = 10000
$ws.Range("P152").Value = 10000
$ws.Range("R152").Value = 'Provincia de Quillota'
$ws.Range("S152").Value = 5000

# Row 153
$ws.Range("D153").Value = 44244
$ws.Range("L153").Value = 'Primera'
$ws.Range("M153").Value = 75
$ws.Range("N153").Value = 3800
$ws.Range("O153").Value = 4000
$ws.Range("P153").Value = 3920
$ws.Range("S153").Value = 1960

# Row 154
$ws.Range("D154").Value = 44179
$ws.Range("M154").Value = 65
$ws.Range("N154").Value = 6000
$ws.Range("O154").Value = 6000
$ws.Range("P154").Value = 6000
$ws.Range("R154").Value = 'Provincia de Curicó'
$ws.Range("S154").Value = 3000

# Row 155
$ws.Range("A155").Value = 3
$ws.Range("B155").Value = 'Femacal de La Calera'
$ws.Range("C155").Value = 'Coquimbo'
$ws.Range("D155").Value = 44179
$ws.Range("E155").Value = 5
$ws.Range("F155").Value = 'Fruta'
$ws.Range("G155").Value = 100101
$ws.Range("H155").Value = 'Berries'
$ws.Range("I155").Value = 100101001
$ws.Range("J155").Value = 'Arándano (blue)'
$ws.Range("K155").Value = 'Sin especificar'
$ws.Range("L155").Value = 'Segunda'
$ws.Range("M155").Value = 60
$ws.Range("N155").Value = 5000
$ws.Range("O155").Value = 5000
$ws.Range("P155").Value = 5000
$ws.Range("Q155").Value = '$/bandeja 2 kilos'
$ws.Range("R155").Value = 'Provincia de Curicó'
$ws.Range("S155").Value = 2500
$ws.Range("T155").Value = 2

# Row 156
$ws.Range("A156").Value = 3
$ws.Range("B156").Value = 'Femacal de La Calera'
$ws.Range("C156").Value = 'Coquimbo'
$ws.Range("D156").Value = 44491
$ws.Range("E156").Value = 5
$ws.Range("F156").Value = 'Fruta'
$ws.Range("G156").Value = 100101
$ws.Range("H156").Value = 'Berries'
$ws.Range("I156").Value = 100101001
$ws.Range("J156").Value = 'Arándano (blue)'
$ws.Range("K156").Value = 'Sin especificar'
$ws.Range("L156").Value = 'Primera'
$ws.Range("M156").Value = 45
$ws.Range("N156").Value = 10000
$ws.Range("O156").Value = 10000
$ws.Range("P156").Value = 10000
$ws.Range("Q156").Value = '$/bandeja 2 kilos'
$ws.Range("R156").Value = 'Provincia de Quillota'
$ws.Range("S156").Value = 5000
$ws.Range("T156").Value = 2

# Apply the date number format used by column D to the newly added rows
$ws.Range("D155").NumberFormat = $ws.Range("D148").NumberFormat()
$ws.Range("D156").NumberFormat = $ws.Range("D148").NumberFormat()
